# Regenerate save_data: column G ("K" - strikeouts) values were recomputed
# (using K instead of Strike#) and rewritten for each saved-game row.
# Only column G changes; all other columns are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G ("K")
$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 3
    16 = 1
    17 = 1
    18 = 0
    19 = 3
    20 = 1
    21 = 1
    22 = 2
    23 = 3
    24 = 0
    25 = 0
    26 = 1
    27 = 3
    28 = 2
    29 = 1
    30 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
